# "Delete description from each excel sheet"
#
# The "device" sheet already carries the real title / update-date / comment
# text in A1:A3. The "Note" sheet only had the generic placeholder labels
# ("Title:", "Update Date:", "Comment:") in A1:A3 - this replaces those
# placeholders with the actual description text (copied from "device"),
# which makes the old generic labels unused and they fall out of the
# shared-string table on save.

$wb     = $excel.ActiveWorkbook
$note   = $wb.Worksheets.Item("Note")
$device = $wb.Worksheets.Item("device")

$title      = $device.Range("A1").Value2
$updateDate = $device.Range("A2").Value2
$comment    = $device.Range("A3").Value2

$note.Range("A1").Value2 = $title
$note.Range("A2").Value2 = $updateDate
$note.Range("A3").Value2 = $comment

# Match the target page setup (same as the "device" sheet).
$note.PageSetup.PaperSize = 9
$note.PageSetup.Orientation = 1

# The "Note" sheet becomes the active tab/selection, "device" loses it.
$note.Activate()
$note.Range("B5").Select() | Out-Null
